$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.240651607513428
$ws.Range("B1").Value = 2.331485271453857
$ws.Range("C1").Value = 2.417727708816528
$ws.Range("D1").Value = 3.19929051399231
$ws.Range("E1").Value = 2.386465311050415
